# css fixes for parser, contacts page
#
# "On top of all routing there is an App_layout class with simply define..."
# becomes
# "On top of all routing there is an App_layout class which simply define..."
#
# and the document's lone "_GoBack" bookmark (left over from the previous
# edit, previously sitting right after "...it will be transfer[ed]") moves
# to sit right after the newly-typed word "which".
#
# Word keeps only a single "_GoBack" bookmark in a document at a time, so
# re-adding it at the new location automatically removes it from the old
# one for us.

$d = $word.ActiveDocument

# 1) Locate the run that needs editing: "class with simply define".
$target = $d.Content
$target.Find.Execute("class with simply define")
if (-not $target.Find.Found) {
    throw "Could not find 'class with simply define' in the document"
}
$phraseStart = $target.Start
$phraseEnd = $target.End

# 2) Narrow down to just the word "with" inside that run, to know exactly
#    where it starts.
$withRange = $d.Range($phraseStart, $phraseEnd)
$withRange.Find.Execute("with")
if (-not $withRange.Find.Found) {
    throw "Could not find 'with' inside the target run"
}
$pos1 = $withRange.Start

# 3) Drop two throwaway bookmarks: one right before "with" and one at the
#    end of the original run (right after "define"). This keeps the edit
#    below from re-coalescing this run with its neighbors (the following
#    "s" run, etc. - which must stay exactly as they were) while still
#    letting it freely rewrite everything from "with" through "define".
$d.Bookmarks.Add("zEditStart", $d.Range($pos1, $pos1))
$d.Bookmarks.Add("zEditEnd", $d.Range($phraseEnd, $phraseEnd))

# 4) Rewrite "with simply define" -> "which simply define" within that
#    isolated range (only the word actually changes).
$tailRange = $d.Range($pos1, $phraseEnd)
$tailRange.Text = "which simply define"

# 5) Remove the throwaway bookmarks again.
$d.Bookmarks("zEditStart").Delete()
$d.Bookmarks("zEditEnd").Delete()

# 6) Re-find the (now updated) phrase and the word "which" inside it, so we
#    can drop the "_GoBack" bookmark immediately after it, collapsed (i.e.
#    bookmarkStart/bookmarkEnd back to back). Adding "_GoBack" here
#    automatically removes it from wherever it used to be.
$after = $d.Content
$after.Find.Execute("class which simply define")
if (-not $after.Find.Found) {
    throw "Could not find 'class which simply define' after the edit"
}
$whichRange = $d.Range($after.Start, $after.End)
$whichRange.Find.Execute("which")
if (-not $whichRange.Find.Found) {
    throw "Could not find 'which' after the edit"
}

$goBackPos = $whichRange.End
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

Write-Output "Replaced 'with' -> 'which' and relocated the _GoBack bookmark"
